$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns retain text formatting so dotted/
# numeric-looking strings (e.g. "30.415.51", "1.005") are not
# auto-converted to numbers when the .Value is assigned below.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.415.51"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "2.104.37"
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("E4").Value = "  +0.99%  "
$ws.Range("D5").Value = "334.57"
$ws.Range("E5").Value = "  +2.15%  "
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("D7").Value = "0.5223"
$ws.Range("E7").Value = "  +0.55%  "
$ws.Range("D8").Value = "0.4561"
$ws.Range("E8").Value = "  +5.58%  "
$ws.Range("D9").Value = "53.38"
$ws.Range("E9").Value = "  +15.85%  "
$ws.Range("D10").Value = "0.08979"
$ws.Range("E10").Value = "  +1.27%  "
$ws.Range("D11").Value = "1.173"
$ws.Range("E11").Value = "  +1.41%  "
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "2.115.53"
$ws.Range("E13").Value = "  +2.15%  "
$ws.Range("D14").Value = "6.837"
$ws.Range("E14").Value = "  +2.61%  "
$ws.Range("D15").Value = "8.023"
$ws.Range("E15").Value = "  +4.55%  "
$ws.Range("D16").Value = "96.48"
$ws.Range("E16").Value = "  +1.53%  "
$ws.Range("D17").Value = "0.00001148"
$ws.Range("E17").Value = "  +2.61%  "
$ws.Range("D18").Value = "1.006"
$ws.Range("D19").Value = "0.06664"
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("D20").Value = "19.21"
$ws.Range("E20").Value = "  +2.42%  "
$ws.Range("D21").Value = "1.005"
$ws.Range("E21").Value = "  +0.69%  "
$ws.Range("D22").Value = "6.340"
$ws.Range("E22").Value = "  +1.90%  "
$ws.Range("D23").Value = "30.490.79"
$ws.Range("E23").Value = "  +0.57%  "
$ws.Range("D24").Value = "12.46"
$ws.Range("E24").Value = "  +1.74%  "
$ws.Range("D25").Value = "2.363"
$ws.Range("E25").Value = "  +3.47%  "
$ws.Range("D26").Value = "2.362.50"
$ws.Range("E26").Value = "  +2.09%  "
$ws.Range("D27").Value = "22.27"
$ws.Range("E27").Value = "  +0.46%  "
$ws.Range("D28").Value = "163.72"
$ws.Range("E28").Value = "  +1.27%  "
$ws.Range("D29").Value = "2.544"
$ws.Range("E29").Value = "  -0.80%  "
$ws.Range("D30").Value = "133.12"
$ws.Range("E30").Value = "  +1.73%  "
$ws.Range("D31").Value = "1.217"
$ws.Range("E31").Value = "  +2.58%  "
$ws.Range("D32").Value = "0.1074"
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("D34").Value = "6.358"
$ws.Range("E34").Value = "  +4.03%  "
$ws.Range("D35").Value = "3.943"
$ws.Range("E35").Value = "  +3.27%  "
$ws.Range("D36").Value = "10.51"
$ws.Range("E36").Value = "  +7.72%  "
$ws.Range("D37").Value = "0.02579"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").Value = "5.684"
$ws.Range("D39").Value = "0.06846"
$ws.Range("E39").Value = "  +3.51%  "
$ws.Range("D40").Value = "0.2304"
$ws.Range("E40").Value = "  +2.88%  "
$ws.Range("D41").Value = "12.70"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "0.6882"
$ws.Range("E42").Value = "  +1.28%  "
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").Value = "2.327"
$ws.Range("E44").Value = "  +5.77%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "14.09"
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.6368"
$ws.Range("E46").Value = "  +0.50%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "3.668"
$ws.Range("E47").Value = "  +1.72%  "
$ws.Range("B48").Value = "EOS"
$ws.Range("C48").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D48").Value = "1.250"
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("B49").Value = "WOONetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D49").Value = "0.3438"
$ws.Range("E49").Value = "  +26.17%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.00000000345"
$ws.Range("E50").Value = "  +17.61%  "
$ws.Range("B51").Value = "WEMIXTOKEN"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").Value = "1.206"
$ws.Range("E51").Value = "  +2.12%  "
